$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.416.22"
$ws.Range("E2").Value = "  +4.35%  "
$ws.Range("D3").Value = "1.794.45"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'314.48"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5458"
$ws.Range("E7").Value = "  +3.75%  "
$ws.Range("D8").Value = "'0.3828"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("D9").Value = "'0.07574"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").Value = "'42.42"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").Value = "'1.121"
$ws.Range("E11").Value = "  +3.13%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'21.10"
$ws.Range("E13").Value = "  +3.33%  "
$ws.Range("D14").Value = "'6.179"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").Value = "'7.389"
$ws.Range("E15").Value = "  +6.69%  "
$ws.Range("D16").Value = "1.794.30"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "'91.55"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "'0.00001068"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").Value = "'0.06453"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("D22").Value = "'5.956"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").Value = "28.411.65"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("D24").Value = "'11.40"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "'2.119"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").Value = "'159.48"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "'2.400"
$ws.Range("D29").Value = "2.004.04"
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("D30").Value = "'122.99"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'1.119"
$ws.Range("E31").Value = "  +5.88%  "
$ws.Range("D32").Value = "'0.1024"
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("D33").Value = "'5.735"
$ws.Range("E33").Value = "  +3.16%  "
$ws.Range("D34").Value = "'3.693"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "'0.2311"
$ws.Range("E35").Value = "  +14.33%  "
$ws.Range("D36").Value = "'0.06409"
$ws.Range("E36").Value = "  +7.50%  "
$ws.Range("D37").Value = "'0.02319"
$ws.Range("E37").Value = "  +3.77%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.151"
$ws.Range("E38").Value = "  +6.72%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'8.771"
$ws.Range("E39").Value = "  +8.98%  "
$ws.Range("D40").Value = "'11.61"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("D41").Value = "'0.6386"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'1.157"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").Value = "'1.389"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "'13.60"
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("D46").Value = "'0.5960"
$ws.Range("D47").Value = "'3.673"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "'125.98"
$ws.Range("E48").Value = "  +3.91%  "
$ws.Range("D49").Value = "'1.986"
$ws.Range("E49").Value = "  +5.95%  "
$ws.Range("D50").Value = "'1.148"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").Value = "'0.06894"
$ws.Range("E51").Value = "  +2.91%  "
